$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "EQUITAS SHELTER" block of rows (rows 6-9), shifting the
# "DRYS SHELTER PREV" rows below it up into their place.
$ws.Rows("6:9").Delete()

$ws.Rows("6:9").Select()
